$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column at J (shifts old J->K, K->L, L->M) ---
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(10).ColumnWidth = 11.6

# --- 2. New column J header + per-row "LoadRunner" Yes/No values ---
$ws.Range("J3").Value = "LoadRunner"
$ws.Range("J3").Style = $ws.Range("K3").Style

$ws.Range("J4").Value = "No"
$ws.Range("J5").Value = "No"
$ws.Range("J6").Value = "No"
$ws.Range("J7").Value = "No"
$ws.Range("J8").Value = "No"
$ws.Range("J9").Value = "No"
$ws.Range("J10").Value = "No"
$ws.Range("J11").Value = "No"
$ws.Range("J12").Value = "No"
$ws.Range("J13").Value = "No"

$ws.Range("J4:J13").Style = $ws.Range("I4:I13").Style

# --- 3. Append new rows 14-17 (copy formatting from row 13 first) ---
$ws.Range("B13:M13").Copy()
$ws.Range("B14:M17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 14
$ws.Range("B14").Value = 44095
$ws.Range("C14").Value = 0.44791666666666669
$ws.Range("D14").Value = 0.625
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = "Off"
$ws.Range("G14").Value = "Yes"
$ws.Range("H14").Value = "2hrs"
$ws.Range("I14").Value = "Yes"
$ws.Range("J14").Value = "No"
$ws.Range("K14").Value = "Success"
$ws.Range("L14").Value = "<70"
$ws.Range("M14").Value = "40+"

# Row 15
$ws.Range("B15").Value = 44095
$ws.Range("C15").Value = 0.63541666666666663
$ws.Range("D15").Value = 0.79861111111111116
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = "On"
$ws.Range("G15").Value = "NA"
$ws.Range("H15").Value = "2hrs"
$ws.Range("I15").Value = "Yes"
$ws.Range("J15").Value = "No"
$ws.Range("K15").Value = "Success"
$ws.Range("L15").Value = "70+"
$ws.Range("M15").Value = 45

# Row 16
$ws.Range("B16").Value = 44097
$ws.Range("C16").Value = 0.53472222222222221
$ws.Range("D16").Value = 0.6875
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = "On"
$ws.Range("G16").Value = "NA"
$ws.Range("H16").Value = "2hrs"
$ws.Range("I16").Value = "Yes"
$ws.Range("J16").Value = "Yes"
$ws.Range("K16").Value = "Success"
$ws.Range("L16").Value = "60+"
$ws.Range("M16").Value = "30+"

# Row 17
$ws.Range("B17").Value = 44097
$ws.Range("C17").Value = 0.69791666666666663
$ws.Range("D17").Value = 0.89583333333333337
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = "On"
$ws.Range("G17").Value = "NA"
$ws.Range("H17").Value = "2hrs"
$ws.Range("I17").Value = "Yes"
$ws.Range("J17").Value = "Yes"
$ws.Range("K17").Value = "Success"
$ws.Range("L17").Value = 80
$ws.Range("M17").Value = "<60"

# --- 4. Final selection, matching the saved view state ---
$ws.Range("N4").Select()
